$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hotel Data")

# Update the Check-In Date and Duration for room 301 (Elon Musk), row 12.
# Force the cell to Text format first so the literal "2024-11-10" string is
# not auto-converted into a date serial number, then clear the formatting
# change back off so the cell keeps its original (default) style.
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "2024-11-10"
$ws.Range("F12").ClearFormats()

$ws.Range("G12").Value = 20
